# Apply the "Enhanced debugging and improved amount extraction" edit:
# The KRA_Database sheet's rows 3-11 are rotated up by one (row 3 takes what
# was row 4's data, row 4 takes what was row 5's data, ... row 10 takes what
# was row 11's data) with the record_id (K) carried along with the moved
# row, the Best_Score (N) normalized to 100 for rows that have real merge
# data, and the "year" (F) column coerced from text to a real number.
# Row 11 ends up holding what used to be row 3's data (the migration test
# record), unchanged. Row 2 (James Mutoro Kitui) is untouched.
# The Database_Summary sheet's "Last Updated" timestamp and "New Records
# Added (This Session)" counter are also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KRA_Database")

# --- Row 3 (was row 4's data: Daisy Jepkosgei Biwott) ---
$ws.Cells.Item(3,1).Value  = "04TH September, 2025"
$ws.Cells.Item(3,2).Value  = "A018905312S"
$ws.Cells.Item(3,3).Value  = "Daisy Jepkosgei Biwott"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value  = "2025."
$ws.Cells.Item(3,6).Value  = 2024
$ws.Cells.Item(3,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(3,8).Value  = "KITALE"
$ws.Cells.Item(3,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(3,10).Value = "multi_format_extractor"
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 2
$ws.Cells.Item(3,13).Value = "Unknown, Unknown"
$ws.Cells.Item(3,14).Value = 100

# --- Row 4 (was row 5's data: Ezekiel Kipserem Korir) ---
$ws.Cells.Item(4,1).Value  = "4th September, 2025"
$ws.Cells.Item(4,2).Value  = "A009775891W"
$ws.Cells.Item(4,3).Value  = "Ezekiel Kipserem Korir"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value  = "2025."
$ws.Cells.Item(4,6).Value  = 2024
$ws.Cells.Item(4,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(4,8).Value  = "KITALE"
$ws.Cells.Item(4,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(4,10).Value = "multi_format_extractor"
$ws.Cells.Item(4,11).Value = 4
$ws.Cells.Item(4,12).Value = 2
$ws.Cells.Item(4,13).Value = "Unknown, Unknown"
$ws.Cells.Item(4,14).Value = 100

# --- Row 5 (was row 6's data: JESSY KAGONDU WAMBUGU) ---
$ws.Cells.Item(5,1).Value  = "04th September, 2025"
$ws.Cells.Item(5,2).Value  = "A004578892U"
$ws.Cells.Item(5,3).Value  = "JESSY KAGONDU WAMBUGU"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value  = "2025."
$ws.Cells.Item(5,6).Value  = 2024
$ws.Cells.Item(5,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(5,8).Value  = "KITALE"
$ws.Cells.Item(5,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(5,10).Value = "multi_format_extractor"
$ws.Cells.Item(5,11).Value = 6
$ws.Cells.Item(5,12).Value = 2
$ws.Cells.Item(5,13).Value = "Unknown, Unknown"
$ws.Cells.Item(5,14).Value = 100

# --- Row 6 (was row 7's data: KELVIN KIPKEMBOI MUTAI) ---
$ws.Cells.Item(6,1).Value  = "04th September, 2025"
$ws.Cells.Item(6,2).Value  = "A008596925K"
$ws.Cells.Item(6,3).Value  = "KELVIN KIPKEMBOI MUTAI"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value  = "2025."
$ws.Cells.Item(6,6).Value  = 2024
$ws.Cells.Item(6,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(6,8).Value  = "KITALE"
$ws.Cells.Item(6,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(6,10).Value = "multi_format_extractor"
$ws.Cells.Item(6,11).Value = 7
$ws.Cells.Item(6,12).Value = 2
$ws.Cells.Item(6,13).Value = "Unknown, Unknown"
$ws.Cells.Item(6,14).Value = 100

# --- Row 7 (was row 8's data: MICHAEL MWANGI MUCHUNGI) ---
$ws.Cells.Item(7,1).Value  = "10th September, 2025"
$ws.Cells.Item(7,2).Value  = "A007388222W"
$ws.Cells.Item(7,3).Value  = "MICHAEL MWANGI MUCHUNGI"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value  = "2025."
$ws.Cells.Item(7,6).Value  = 2024
$ws.Cells.Item(7,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(7,8).Value  = "KITALE"
$ws.Cells.Item(7,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(7,10).Value = "multi_format_extractor"
$ws.Cells.Item(7,11).Value = 8
$ws.Cells.Item(7,12).Value = 2
$ws.Cells.Item(7,13).Value = "Unknown, Unknown"
$ws.Cells.Item(7,14).Value = 100

# --- Row 8 (was row 9's data: Paul Chotomolo Mirikwa) ---
$ws.Cells.Item(8,1).Value  = "04th September, 2025"
$ws.Cells.Item(8,2).Value  = "A012209532N"
$ws.Cells.Item(8,3).Value  = "Paul Chotomolo Mirikwa"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value  = "2025."
$ws.Cells.Item(8,6).Value  = 2024
$ws.Cells.Item(8,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(8,8).Value  = "NAITIRI"
$ws.Cells.Item(8,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(8,10).Value = "multi_format_extractor"
$ws.Cells.Item(8,11).Value = 9
$ws.Cells.Item(8,12).Value = 2
$ws.Cells.Item(8,13).Value = "Unknown, Unknown"
$ws.Cells.Item(8,14).Value = 100

# --- Row 9 (was row 10's data: Peter Kimutai Telengech) ---
$ws.Cells.Item(9,1).Value  = "29TH AUGUST, 2025"
$ws.Cells.Item(9,2).Value  = "A001126762Z"
$ws.Cells.Item(9,3).Value  = "Peter Kimutai Telengech"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value  = "2025."
$ws.Cells.Item(9,6).Value  = 2024
$ws.Cells.Item(9,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(9,8).Value  = "ELDORET"
$ws.Cells.Item(9,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(9,10).Value = "multi_format_extractor"
$ws.Cells.Item(9,11).Value = 10
$ws.Cells.Item(9,12).Value = 2
$ws.Cells.Item(9,13).Value = "Unknown, Unknown"
$ws.Cells.Item(9,14).Value = 100

# --- Row 10 (was row 11's data: THOMAS JUMA SIKUKU) ---
$ws.Cells.Item(10,1).Value  = "04th September, 2025"
$ws.Cells.Item(10,2).Value  = "A005615142S"
$ws.Cells.Item(10,3).Value  = "THOMAS JUMA SIKUKU"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value  = "2025"
$ws.Cells.Item(10,6).Value  = 2024
$ws.Cells.Item(10,7).Value  = "Franciscar Nyangweta"
$ws.Cells.Item(10,8).Value  = "KITALE"
$ws.Cells.Item(10,9).Value  = "2025-09-22 11:25:52"
$ws.Cells.Item(10,10).Value = "multi_format_extractor"
$ws.Cells.Item(10,11).Value = 11
$ws.Cells.Item(10,12).Value = 2
$ws.Cells.Item(10,13).Value = "Unknown, Unknown"
$ws.Cells.Item(10,14).Value = 100

# --- Row 11 (was row 3's data: Test User After Migration) ---
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value  = "2024-09-22"
$ws.Cells.Item(11,2).Value  = "A123456789X"
$ws.Cells.Item(11,3).Value  = "Test User After Migration"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value  = "15,000.00"
$ws.Cells.Item(11,6).Value  = 2024
$ws.Cells.Item(11,7).Value  = "Test Officer"
$ws.Cells.Item(11,8).Value  = "NAIROBI"
$ws.Cells.Item(11,9).Value  = "2025-09-22 09:11:27"
$ws.Cells.Item(11,10).Value = "post_migration_test"
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).ClearContents()
$ws.Cells.Item(11,13).ClearContents()
$ws.Cells.Item(11,14).ClearContents()

# --- Database_Summary sheet: refresh the "Last Updated" timestamp and the
#     "New Records Added (This Session)" counter ---
$summary = $wb.Worksheets.Item("Database_Summary")
$summary.Cells.Item(3,2).Value = "2025-09-22 12:16:12"
$summary.Cells.Item(4,2).Value = 1
